$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row height for existing row 53 (13.8 -> 12.8, matching the rest of the sheet)
$ws.Rows.Item(53).RowHeight = 12.8

# New KPI rows appended after row 53
$newRows = @(
    @{ Row = 54; Id = 4601; Text = "Нет ЦА: Категория тов. для животных примыкает или расположена в радиусе 5 м от центра выкладки приоритетной категории (1) (молочные прод, фрукты и овощи, хлебобулочные изд, кондитерские изд, мясн. изд. и рыба) таким образом, что видны блоки паучей Kitekat и Whiskas" },
    @{ Row = 55; Id = 4602; Text = "Нет ЦА: Категория товаров для животных примыкает или расположена в радиусе 5 м от центра выкладки приоритетной категории (2) (консервы, соки, вода/газированные напитки, замороженные продукты), таким образом, что видны блоки паучей Kitekat  и Whiskas" },
    @{ Row = 56; Id = 4603; Text = "Нет ЦА: Категория выстроена в единую линию единым блоком или образует внутренний угол" },
    @{ Row = 57; Id = 4604; Text = "Нет ЦА: Категория располагается вне тупика и находится дальше 5 м от входа/выхода и кассовой зоны" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = "MARS KPIs"
    $ws.Cells.Item($row, 2).Value = $r.Id
    $ws.Cells.Item($row, 3).Value = $r.Id
    $ws.Cells.Item($row, 4).Value = $r.Text
    $ws.Cells.Item($row, 2).HorizontalAlignment = -4131
    $ws.Cells.Item($row, 3).HorizontalAlignment = -4131
    $ws.Rows.Item($row).RowHeight = 12.8
}

# Update selection/view to match post-edit state
$ws.Range("A55").Select()
$excel.ActiveWindow.ScrollRow = 31
